# Auto-generated script to apply scheduled-runner data refresh to Seraph_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 733.6667
$ws.Range("I2").Value = 462.25
$ws.Range("J2").Value = 1276.5
$ws.Range("K2").Value = 462.25
$ws.Range("L2").Value = 1276.5
$ws.Range("M2").Value = -349.25
$ws.Range("N2").Value = -1502.5
$ws.Range("H28").Value = 1292.7
$ws.Range("I28").Value = 991
$ws.Range("K28").Value = 991
$ws.Range("M28").Value = -506
$ws.Range("H40").Value = 2350
$ws.Range("I40").Value = 1833.3334
$ws.Range("J40").Value = 2571.4285
$ws.Range("K40").Value = 1833.3334
$ws.Range("L40").Value = 2571.4285
$ws.Range("M40").Value = -1658.3334
$ws.Range("N40").Value = -2921.4285
$ws.Range("H51").Value = 5000
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -5968
$ws.Range("H86").Value = 3612.3333
$ws.Range("I86").Value = 2993.75
$ws.Range("K86").Value = 2993.75
$ws.Range("M86").Value = -1870.75
$ws.Range("H89").Value = 3612.3333
$ws.Range("I89").Value = 2993.75
$ws.Range("K89").Value = 14968.75
$ws.Range("M89").Value = -9352.75
$ws.Range("H106").Value = 19751.166
$ws.Range("I106").Value = 22101.4
$ws.Range("K106").Value = 22101.4
$ws.Range("M106").Value = -21470.4
$ws.Range("H116").Value = 4785
$ws.Range("I116").Value = 4785
$ws.Range("K116").Value = 4785
$ws.Range("M116").Value = -1343
$ws.Range("H132").Value = 1463.2
$ws.Range("I132").Value = 1462.2307
$ws.Range("K132").Value = 4386.6921
$ws.Range("M132").Value = -1856.6921
$ws.Range("H137").Value = 2277.4546
$ws.Range("I137").Value = 2464.7144
$ws.Range("K137").Value = 7394.1432
$ws.Range("M137").Value = -4844.1432
$ws.Range("H138").Value = 6750.1562
$ws.Range("I138").Value = 890.3889
$ws.Range("J138").Value = 14284.143
$ws.Range("K138").Value = 2671.1667
$ws.Range("L138").Value = 42852.429
$ws.Range("M138").Value = 2468.8333
$ws.Range("N138").Value = -53132.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4164.5
$ws.Range("I32").Value = 3166.081
$ws.Range("K32").Value = 3166.081
$ws.Range("M32").Value = -2879.081
$ws.Range("H61").Value = 2980.8572
$ws.Range("I61").Value = 3378.2
$ws.Range("K61").Value = 3378.2
$ws.Range("M61").Value = -3166.2
$ws.Range("H132").Value = 1300.5454
$ws.Range("I132").Value = 1320.6
$ws.Range("K132").Value = 3961.8
$ws.Range("M132").Value = -1431.8
$ws.Range("H136").Value = 2980.8572
$ws.Range("I136").Value = 3378.2
$ws.Range("K136").Value = 10134.6
$ws.Range("M136").Value = -7584.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4791.375
$ws.Range("I86").Value = 3566.6667
$ws.Range("J86").Value = 5526.2
$ws.Range("K86").Value = 3566.6667
$ws.Range("L86").Value = 5526.2
$ws.Range("M86").Value = -2443.6667
$ws.Range("N86").Value = -7772.2
$ws.Range("H89").Value = 4791.375
$ws.Range("I89").Value = 3566.6667
$ws.Range("J89").Value = 5526.2
$ws.Range("K89").Value = 17833.3335
$ws.Range("L89").Value = 27631
$ws.Range("M89").Value = -12217.3335
$ws.Range("N89").Value = -38863
$ws.Range("H105").Value = 2677.5
$ws.Range("I105").Value = 2744
$ws.Range("J105").Value = 2345
$ws.Range("K105").Value = 2744
$ws.Range("L105").Value = 2345
$ws.Range("M105").Value = -997
$ws.Range("N105").Value = -5839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 970.46155
$ws.Range("J107").Value = 1316.3636
$ws.Range("L107").Value = 1316.3636
$ws.Range("N107").Value = -5156.3636
$ws.Range("H122").Value = 3033.4707
$ws.Range("I122").Value = 3214.1667
$ws.Range("J122").Value = 2599.8
$ws.Range("K122").Value = 9642.500100000001
$ws.Range("L122").Value = 7799.400000000001
$ws.Range("M122").Value = -7192.500100000001
$ws.Range("N122").Value = -12699.4
$ws.Range("H132").Value = 1982
$ws.Range("I132").Value = 1548.7028
$ws.Range("K132").Value = 4646.1084
$ws.Range("M132").Value = -2116.1084
$ws.Range("H134").Value = 1752.3773
$ws.Range("J134").Value = 2498
$ws.Range("L134").Value = 7494
$ws.Range("N134").Value = -12564

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1039751.44
$ws.Range("J4").Value = 715.375
$ws.Range("L4").Value = 2146.125
$ws.Range("N4").Value = -2370.125
$ws.Range("H68").Value = 634
$ws.Range("I68").Value = 451
$ws.Range("K68").Value = 1353
$ws.Range("M68").Value = -542
$ws.Range("H71").Value = 634
$ws.Range("I71").Value = 451
$ws.Range("K71").Value = 4059
$ws.Range("M71").Value = -3
$ws.Range("H132").Value = 9254.933999999999
$ws.Range("J132").Value = 4924.5
$ws.Range("L132").Value = 44320.5
$ws.Range("N132").Value = -49380.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2913.4736
$ws.Range("J132").Value = 2925.3809
$ws.Range("L132").Value = 8776.1427
$ws.Range("N132").Value = -13836.1427

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 19566.334
$ws.Range("J16").Value = 19350
$ws.Range("L16").Value = 19350
$ws.Range("N16").Value = -19690
$ws.Range("H46").Value = 3079.9473
$ws.Range("I46").Value = 2232.2307
$ws.Range("J46").Value = 4916.6665
$ws.Range("K46").Value = 2232.2307
$ws.Range("L46").Value = 4916.6665
$ws.Range("M46").Value = -2044.2307
$ws.Range("N46").Value = -5292.6665
$ws.Range("H82").Value = 1792.1428
$ws.Range("I82").Value = 1436.25
$ws.Range("J82").Value = 2266.6667
$ws.Range("K82").Value = 1436.25
$ws.Range("L82").Value = 2266.6667
$ws.Range("M82").Value = -1075.25
$ws.Range("N82").Value = -2988.6667
$ws.Range("H85").Value = 1792.1428
$ws.Range("I85").Value = 1436.25
$ws.Range("J85").Value = 2266.6667
$ws.Range("K85").Value = 1436.25
$ws.Range("L85").Value = 2266.6667
$ws.Range("M85").Value = -188.25
$ws.Range("N85").Value = -4762.6667
$ws.Range("H100").Value = 3749
$ws.Range("J100").Value = 4000
$ws.Range("L100").Value = 4000
$ws.Range("N100").Value = -5082
$ws.Range("H101").Value = 5000
$ws.Range("J101").Value = 5000
$ws.Range("L101").Value = 5000
$ws.Range("N101").Value = -11490
$ws.Range("H122").Value = 5199.8
$ws.Range("I122").Value = 5333.3335
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 16000.0005
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -13550.0005
$ws.Range("N122").Value = -19898.5
$ws.Range("H132").Value = 5131.591
$ws.Range("I132").Value = 4165.385
$ws.Range("K132").Value = 12496.155
$ws.Range("M132").Value = -9966.155000000001
$ws.Range("H136").Value = 2898.8948
$ws.Range("I136").Value = 2699.1428
$ws.Range("K136").Value = 8097.428400000001
$ws.Range("M136").Value = -5547.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 20333.334
$ws.Range("J54").Value = 20333.334
$ws.Range("L54").Value = 20333.334
$ws.Range("N54").Value = -21373.334
$ws.Range("H103").Value = 41375
$ws.Range("J103").Value = 41375
$ws.Range("L103").Value = 41375
$ws.Range("N103").Value = -43719
$ws.Range("H129").Value = 40000
$ws.Range("J129").Value = 40000
$ws.Range("L129").Value = 40000
$ws.Range("N129").Value = -50000
$ws.Range("H132").Value = 34867.32
$ws.Range("I132").Value = 45944.906
$ws.Range("K132").Value = 137834.718
$ws.Range("M132").Value = -135304.718
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
